$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.28"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'0.34%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'40.86"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'0.04%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.127"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'1.65%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.07623"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-0.22%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'1.606"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-0.06%"
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'1.63%"
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'-0.07%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.1267"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'25.96%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1807"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'2.38%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.08993"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-1.90%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.04292"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'1.58%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.1046"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-0.62%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.001253"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-0.05%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.005739"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-1.31%"
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'-0.29%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'4.292"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'0.76%"
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'1.41%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'6.915"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'1.78%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.1385"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'2.13%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.2741"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'0.71%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.04046"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-2.59%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.001271"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'4.18%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.004039"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-0.95%"
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'-2.12%"
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'24.81%"
$ws.Range("E26").ClearFormats()
$ws.Range("D38").Value = "'0.02413"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'-0.37%"
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.05230"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'0.98%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.007842"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'0.76%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.1299"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-0.71%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.006799"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-3.43%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.001863"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-4.33%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.007431"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-0.65%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.3358"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'9.85%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006882"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'8.27%"
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'0.25%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.1594"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'2,912.44%"
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'-31.72%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'0.25%"
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'0.25%"
$ws.Range("E51").ClearFormats()
